# Update "想去人数" (F column) counts for several events in both the
# "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 34
$wsExhibit.Range("F7").Value  = 2630
$wsExhibit.Range("F9").Value  = 231
$wsExhibit.Range("F10").Value = 93
$wsExhibit.Range("F11").Value = 5906
$wsExhibit.Range("F15").Value = 11581
$wsExhibit.Range("F16").Value = 11774
$wsExhibit.Range("F18").Value = 76
$wsExhibit.Range("F21").Value = 61

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 34
$wsAll.Range("F7").Value  = 2630
$wsAll.Range("F10").Value = 231
$wsAll.Range("F11").Value = 93
$wsAll.Range("F12").Value = 5906
$wsAll.Range("F16").Value = 11581
$wsAll.Range("F17").Value = 11774
$wsAll.Range("F19").Value = 76
$wsAll.Range("F22").Value = 61
